$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old Table1 (structured table) that lived over B3:F7 ---
if ($ws.ListObjects.Count -gt 0) {
    $ws.ListObjects.Item(1).Delete()
}

# --- Wipe the old demo content (header "some text" + Column1..Column5) ---
$ws.Cells.Clear()

# Drop the now-unused custom widths that used to span columns C:G together
# with B, so only column B keeps an explicit width definition.
$ws.Range("C1:G1").EntireColumn.Delete()

# --- Rename the sheet from "Sheet1" to "Dashboard" ---
$ws.Name = "Dashboard"

# --- Populate the new dashboard "menu" of named ranges / dataframes ---
# Column B first (top to bottom), then column C (top to bottom) - this
# mirrors how the sheet was actually authored.
$ws.Range("B5").Value = "var_structured_position_top10"
$ws.Range("B6").Value = "fund_exp_pct_dashboard"
$ws.Range("B7").Value = "sector_exposure_df"
$ws.Range("C5").Value = "var_structured_position_bottom10"
$ws.Range("C6").Value = "fund_exp_usd_dashboard"
$ws.Range("B8").Value = "something_else"
$ws.Range("B9").Value = "another thing"
$ws.Range("C9").Value = "another thing to the left"

# --- Column B is widened to fit the longest label ---
$ws.Columns.Item(2).ColumnWidth = 28.7109375

# --- Leave the selection where the author left it ---
$ws.Range("P7").Select() | Out-Null
